$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new "Save" column (H), matching the style already used by the
# other header cells (B1:G1, style index 1: bold, centered, bordered).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Save values for rows 2-14 (0/1 flags)
$saveValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 1
    6  = 1
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 1
    13 = 0
    14 = 1
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}

$excel.CutCopyMode = $false
